$d = $word.ActiveDocument

# Merge the word-by-word runs in the Title, Author and Abstract paragraphs
# into single runs containing the full text, by performing an in-place
# Find/Replace of the whole (already-correct) text. Word's Find/Replace
# re-writes the matched range as a single run, collapsing the many
# single-word runs the document currently has.

$wdReplaceAll = 2
$wdFindContinue = 1

$d.Content.Find.Execute(
    "Answers: Trigonometric identities (radians)", $true, $false, $false,
    $false, $false, $true, $wdFindContinue, $false,
    "Answers: Trigonometric identities (radians)", $wdReplaceAll)

$d.Content.Find.Execute(
    "Dzhemma Ruseva", $true, $false, $false,
    $false, $false, $true, $wdFindContinue, $false,
    "Dzhemma Ruseva", $wdReplaceAll)

$d.Content.Find.Execute(
    "A selection of questions on trigonometric identities, using radians to measure angles.",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "A selection of questions on trigonometric identities, using radians to measure angles.",
    $wdReplaceAll)
